$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove obsolete columns H:L (theta, conformation, monosaccharides, motifs, class)
$ws.Range("H1:L9").EntireColumn.Delete()

# Update header row (A1, B1, F1 stay the same)
$ws.Cells.Item(1,3).Value = "monosaccharides"
$ws.Cells.Item(1,4).Value = "motifs"
$ws.Cells.Item(1,5).Value = "sasa"
$ws.Cells.Item(1,6).Value = "flexibility"
$ws.Cells.Item(1,7).Value = "has_multi_node_motifs"

# Row 2 - update with new per-motif values (glycan/binding_score in A2/B2 unchanged)
$ws.Cells.Item(2,3).Value = "['Fuc(a1-2)', 'Gal(a1-1)', 'Gal(a1-3)']"
$ws.Cells.Item(2,4).Value = "['Fuc(a1-2)[Gal(a1-3)]Gal']"
$ws.Cells.Item(2,5).Value = 6.385619566891085
$ws.Cells.Item(2,6).Value = 1.404918654625511
$ws.Cells.Item(2,7).Value = $true

# Row 3 - new row
$ws.Cells.Item(3,1).Value = "Fuc(a1-2)[Gal(a1-3)]Gal(b1-3)GalNAc"
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(3,1).PasteSpecial(-4122)
$ws.Cells.Item(3,2).Value = 0.9078085701541632
$ws.Cells.Item(3,3).Value = "['Fuc(a1-2)', 'Gal(b1-3)', 'Gal(a1-3)']"
$ws.Cells.Item(3,4).Value = "['Fuc(a1-2)[Gal(a1-3)]Gal']"
$ws.Cells.Item(3,5).Value = 5.727669409197825
$ws.Cells.Item(3,6).Value = 1.532193243723363
$ws.Cells.Item(3,7).Value = $true

# Row 4 - new row
$ws.Cells.Item(4,1).Value = "Fuc(a1-2)[Gal(a1-3)]Gal(b1-3)GalNAc(b1-3)Gal(a1-4)Gal(b1-4)Glc"
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(4,1).PasteSpecial(-4122)
$ws.Cells.Item(4,2).Value = 0.8228596363115226
$ws.Cells.Item(4,3).Value = "['Fuc(a1-2)', 'Gal(b1-3)', 'Gal(a1-3)']"
$ws.Cells.Item(4,4).Value = "['Fuc(a1-2)[Gal(a1-3)]Gal']"
$ws.Cells.Item(4,5).Value = 5.444661550637573
$ws.Cells.Item(4,6).Value = 7.120285194243348
$ws.Cells.Item(4,7).Value = $true

# Row 5 - new row
$ws.Cells.Item(5,1).Value = "Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)Glc"
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(5,1).PasteSpecial(-4122)
$ws.Cells.Item(5,2).Value = 3.094286373856975
$ws.Cells.Item(5,3).Value = "['Fuc(a1-2)', 'Gal(b1-4)', 'Gal(a1-3)']"
$ws.Cells.Item(5,4).Value = "['Fuc(a1-2)[Gal(a1-3)]Gal']"
$ws.Cells.Item(5,5).Value = 5.841043650482035
$ws.Cells.Item(5,6).Value = 1.754103373790929
$ws.Cells.Item(5,7).Value = $true

# Row 6 - new row
$ws.Cells.Item(6,1).Value = "Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc"
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(6,1).PasteSpecial(-4122)
$ws.Cells.Item(6,2).Value = -0.2298719333647885
$ws.Cells.Item(6,3).Value = "['Fuc(a1-2)', 'Gal(b1-4)', 'Gal(a1-3)']"
$ws.Cells.Item(6,4).Value = "['Fuc(a1-2)[Gal(a1-3)]Gal']"
$ws.Cells.Item(6,5).Value = 5.948417738735455
$ws.Cells.Item(6,6).Value = 2.055811179438615
$ws.Cells.Item(6,7).Value = $true

# Row 7 - new row
$ws.Cells.Item(7,1).Value = "Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)GlcNAc(b1-3)GalNAc"
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(7,1).PasteSpecial(-4122)
$ws.Cells.Item(7,2).Value = 0.0122707804677351
$ws.Cells.Item(7,3).Value = "['Fuc(a1-2)', 'Gal(b1-4)', 'Gal(a1-3)']"
$ws.Cells.Item(7,4).Value = "['Fuc(a1-2)[Gal(a1-3)]Gal']"
$ws.Cells.Item(7,5).Value = 5.811879060637235
$ws.Cells.Item(7,6).Value = 1.558149379507167
$ws.Cells.Item(7,7).Value = $true

# Row 8 - new row
$ws.Cells.Item(8,1).Value = "Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc"
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(8,1).PasteSpecial(-4122)
$ws.Cells.Item(8,2).Value = -0.4964489684821903
$ws.Cells.Item(8,3).Value = "['Gal(b1-4)', 'GlcNAc(b1-3)', 'Gal(b1-4)', 'GlcNAc(b1-3)']"
$ws.Cells.Item(8,4).Value = "['Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc(b1-3)']"
$ws.Cells.Item(8,5).Value = 9.740877345841087
$ws.Cells.Item(8,6).Value = 4.372965213239025
$ws.Cells.Item(8,7).Value = $true

# Row 9 - new row
$ws.Cells.Item(9,1).Value = "Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc"
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(9,1).PasteSpecial(-4122)
$ws.Cells.Item(9,2).Value = -0.1227473274964418
$ws.Cells.Item(9,3).Value = "['Gal(b1-4)', 'GlcNAc(b1-3)', 'Gal(b1-4)', 'GlcNAc(b1-3)']"
$ws.Cells.Item(9,4).Value = "['Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc(b1-3)']"
$ws.Cells.Item(9,5).Value = 9.738054376039319
$ws.Cells.Item(9,6).Value = 3.272015281964459
$ws.Cells.Item(9,7).Value = $true
